$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F5").Value = 2
$ws.Range("F10").Value = -1
$ws.Range("F14").Value = -2
$ws.Range("F15").Value = 2
$ws.Range("F17").Value = 2
$ws.Range("F19").Value = -4
$ws.Range("F21").Value = 1
$ws.Range("F22").Value = -1
$ws.Range("F25").Value = -2
$ws.Range("F35").Value = 2
$ws.Range("F36").Value = 0
$ws.Range("F37").Value = 0
$ws.Range("F40").Value = -2
$ws.Range("F48").Value = 4
$ws.Range("F49").Value = -5
$ws.Range("F60").Value = -1
